$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.066889322219749
$ws.Range("D2").Value = 1.063090312354158
$ws.Range("E2").Value = 1.07049453537711
$ws.Range("F2").Value = 1.077910778299346
$ws.Range("I2").Value = 1.050874459765621
$ws.Range("J2").Value = 1.071837020696116
$ws.Range("K2").Value = 1.06580989298125
$ws.Range("L2").Value = 1.073194222558389
$ws.Range("M2").Value = 1.080590835772935
$ws.Range("N2").Value = 1.073359152206929

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.068195889829631
$ws.Range("D3").Value = 1.063897921716811
$ws.Range("E3").Value = 1.07161242045496
$ws.Range("F3").Value = 1.079001568284955
$ws.Range("I3").Value = 1.051209767852282
$ws.Range("J3").Value = 1.07279806718227
$ws.Range("K3").Value = 1.066432463854492
$ws.Range("L3").Value = 1.074127734399043
$ws.Range("M3").Value = 1.08149874225235
$ws.Range("N3").Value = 1.074321563489327

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.06904113068135
$ws.Range("D4").Value = 1.064420370448419
$ws.Range("E4").Value = 1.072335862371256
$ws.Range("F4").Value = 1.079707549650958
$ws.Range("I4").Value = 1.051425447990017
$ws.Range("J4").Value = 1.073419211330042
$ws.Range("K4").Value = 1.066834525618838
$ws.Range("L4").Value = 1.074731281889302
$ws.Range("M4").Value = 1.082085785593998
$ws.Range("N4").Value = 1.074943589733035

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.069396425574413
$ws.Range("D5").Value = 1.064639977293519
$ws.Range("E5").Value = 1.07264002149028
$ws.Range("F5").Value = 1.080004385206392
$ws.Range("I5").Value = 1.05151581231297
$ws.Range("J5").Value = 1.07368017010752
$ws.Range("K5").Value = 1.067003365295397
$ws.Range("L5").Value = 1.074984895454406
$ws.Range("M5").Value = 1.082332475906763
$ws.Range("N5").Value = 1.07520491910191

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.069456078621787
$ws.Range("D6").Value = 1.06467684845946
$ws.Range("E6").Value = 1.072691092597042
$ws.Range("F6").Value = 1.080054227633078
$ws.Range("I6").Value = 1.051530966854157
$ws.Range("J6").Value = 1.073723976295479
$ws.Range("K6").Value = 1.067031703261061
$ws.Range("L6").Value = 1.075027471415113
$ws.Range("M6").Value = 1.082373890299615
$ws.Range("N6").Value = 1.07524878749968

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.069045878320275
$ws.Range("D7").Value = 1.064423304965874
$ws.Range("E7").Value = 1.072339926463862
$ws.Range("F7").Value = 1.079711515819952
$ws.Range("I7").Value = 1.051426656650694
$ws.Range("J7").Value = 1.073422698941017
$ws.Range("K7").Value = 1.066836782397478
$ws.Range("L7").Value = 1.074734671147543
$ws.Range("M7").Value = 1.082089082284458
$ws.Range("N7").Value = 1.074947082296818

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.067330924901528
$ws.Range("D8").Value = 1.063363274435434
$ws.Range("E8").Value = 1.070872310860112
$ws.Range("F8").Value = 1.078279381804461
$ws.Range("I8").Value = 1.050988045353611
$ws.Range("J8").Value = 1.072161960243446
$ws.Range("K8").Value = 1.066020456144844
$ws.Range("L8").Value = 1.073509810908125
$ws.Range("M8").Value = 1.080897757199548
$ws.Range("N8").Value = 1.073684553205688

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.064307346325472
$ws.Range("D9").Value = 1.061494372583932
$ws.Range("E9").Value = 1.068286858939489
$ws.Range("F9").Value = 1.075757017625968
$ws.Range("I9").Value = 1.050205283624735
$ws.Range("J9").Value = 1.069934818388327
$ws.Range("K9").Value = 1.064575967381846
$ws.Range("L9").Value = 1.071347589195151
$ws.Range("M9").Value = 1.078795133356556
$ws.Range("N9").Value = 1.071454248553677

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.06229038478425
$ws.Range("D10").Value = 1.060247756251525
$ws.Range("E10").Value = 1.066563595002596
$ws.Range("F10").Value = 1.07407621235684
$ws.Range("I10").Value = 1.049676777365642
$ws.Range("J10").Value = 1.068446224526607
$ws.Range("K10").Value = 1.063608900149022
$ws.Range("L10").Value = 1.069903433576693
$ws.Range("M10").Value = 1.077391066507114
$ws.Range("N10").Value = 1.06996354071793

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.061416687445361
$ws.Range("D11").Value = 1.059707791817614
$ws.Range("E11").Value = 1.065817471407221
$ws.Range("F11").Value = 1.073348573246944
$ws.Range("I11").Value = 1.049446340948368
$ws.Range("J11").Value = 1.067800716089269
$ws.Range("K11").Value = 1.063189175792693
$ws.Range("L11").Value = 1.069277445774797
$ws.Range("M11").Value = 1.076782526074771
$ws.Range("N11").Value = 1.069317115584569

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.061092103124385
$ws.Range("D12").Value = 1.059507198791269
$ws.Range("E12").Value = 1.065540335304793
$ws.Range("F12").Value = 1.073078318452346
$ws.Range("I12").Value = 1.049360507142256
$ws.Range("J12").Value = 1.067560802696545
$ws.Range("K12").Value = 1.063033123867829
$ws.Range("L12").Value = 1.06904482538817
$ws.Range("M12").Value = 1.076556400067228
$ws.Range("N12").Value = 1.06907686148731

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.061161730068358
$ws.Range("D13").Value = 1.059550227869004
$ws.Range("E13").Value = 1.065599781631858
$ws.Range("F13").Value = 1.073136288014896
$ws.Range("I13").Value = 1.049378929631154
$ws.Range("J13").Value = 1.06761227142824
$ws.Range("K13").Value = 1.063066604226391
$ws.Range("L13").Value = 1.069094727831764
$ws.Range("M13").Value = 1.076604908821546
$ws.Range("N13").Value = 1.069128403310507

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.061389858328867
$ws.Range("D14").Value = 1.059691211272464
$ws.Range("E14").Value = 1.065794563109287
$ws.Range("F14").Value = 1.073326233444776
$ws.Range("I14").Value = 1.04943925078771
$ws.Range("J14").Value = 1.067780887698257
$ws.Range("K14").Value = 1.063176279507168
$ws.Range("L14").Value = 1.069258219363784
$ws.Range("M14").Value = 1.076763836210823
$ws.Range("N14").Value = 1.069297259034968

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.061530408397954
$ws.Range("D15").Value = 1.059778072315814
$ws.Range("E15").Value = 1.065914575337122
$ws.Range("F15").Value = 1.073443268069169
$ws.Range("I15").Value = 1.049476384893603
$ws.Range("J15").Value = 1.06788475876339
$ws.Range("K15").Value = 1.063243834486597
$ws.Range("L15").Value = 1.069358938506565
$ws.Range("M15").Value = 1.076861745053079
$ws.Range("N15").Value = 1.069401277608928

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.062348361769948
$ws.Range("D16").Value = 1.060283588264364
$ws.Range("E16").Value = 1.066613113905437
$ws.Range("F16").Value = 1.074124506636437
$ws.Range("I16").Value = 1.049692037125011
$ws.Range("J16").Value = 1.068489044875897
$ws.Range("K16").Value = 1.063636735210312
$ws.Range("L16").Value = 1.069944964302351
$ws.Range("M16").Value = 1.077431441204434
$ws.Range("N16").Value = 1.070006421877027

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.062861349113779
$ws.Range("D17").Value = 1.060600638874121
$ws.Range("E17").Value = 1.067051303806702
$ws.Range("F17").Value = 1.07455187164071
$ws.Range("I17").Value = 1.049826884019514
$ws.Range("J17").Value = 1.068867845209413
$ws.Range("K17").Value = 1.063882929160666
$ws.Range("L17").Value = 1.070312385231409
$ws.Range("M17").Value = 1.077788642949394
$ws.Range("N17").Value = 1.070385760150466

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.063160533395276
$ws.Range("D18").Value = 1.06078555263135
$ws.Range("E18").Value = 1.067306898573872
$ws.Range("F18").Value = 1.074801162012969
$ws.Range("I18").Value = 1.049905384589472
$ws.Range("J18").Value = 1.069088702717293
$ws.Range("K18").Value = 1.064026435625292
$ws.Range("L18").Value = 1.070526632182865
$ws.Range("M18").Value = 1.077996937756744
$ws.Range("N18").Value = 1.07060693130134

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.063262542077477
$ws.Range("D19").Value = 1.060848600658545
$ws.Range("E19").Value = 1.067394050856369
$ws.Range("F19").Value = 1.0748861662738
$ws.Range("I19").Value = 1.049932125264316
$ws.Range("J19").Value = 1.069163994159734
$ws.Range("K19").Value = 1.064075351646752
$ws.Range("L19").Value = 1.070599674142273
$ws.Range("M19").Value = 1.078067951626258
$ws.Range("N19").Value = 1.070682329666265

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.062806313818614
$ws.Range("D20").Value = 1.06056662405481
$ws.Range("E20").Value = 1.067004289569206
$ws.Range("F20").Value = 1.074506017818031
$ws.Range("I20").Value = 1.049812432089755
$ws.Range("J20").Value = 1.06882721287583
$ws.Range("K20").Value = 1.063856524639455
$ws.Range("L20").Value = 1.070272971019841
$ws.Range("M20").Value = 1.077750324274525
$ws.Range("N20").Value = 1.070345070114308

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.061322681764741
$ws.Range("D21").Value = 1.059649695900804
$ws.Range("E21").Value = 1.0657372046282
$ws.Range("F21").Value = 1.073270298632433
$ws.Range("I21").Value = 1.049421494325865
$ws.Range("J21").Value = 1.067731238367732
$ws.Range("K21").Value = 1.063143986950968
$ws.Range("L21").Value = 1.069210077974128
$ws.Range("M21").Value = 1.07671703846803
$ws.Range("N21").Value = 1.069247539196699

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.060389548074708
$ws.Range("D22").Value = 1.05907303500026
$ws.Range("E22").Value = 1.064940580105136
$ws.Range("F22").Value = 1.072493483986605
$ws.Range("I22").Value = 1.04917431048293
$ws.Range("J22").Value = 1.067041328414593
$ws.Range("K22").Value = 1.06269513184698
$ws.Range("L22").Value = 1.068541211806905
$ws.Range("M22").Value = 1.076066866679416
$ws.Range("N22").Value = 1.068556649492292

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.060884250856571
$ws.Range("D23").Value = 1.059378748393409
$ws.Range("E23").Value = 1.065362882435471
$ws.Range("F23").Value = 1.072905276163982
$ws.Range("I23").Value = 1.049305478876582
$ws.Range("J23").Value = 1.067407141810462
$ws.Range("K23").Value = 1.062933159649137
$ws.Range("L23").Value = 1.068895846243926
$ws.Range("M23").Value = 1.076411583262415
$ws.Range("N23").Value = 1.068922982385144

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.062831182000085
$ws.Range("D24").Value = 1.060581993938643
$ws.Range("E24").Value = 1.067025533260928
$ws.Range("F24").Value = 1.074526737138672
$ws.Range("I24").Value = 1.049818962768579
$ws.Range("J24").Value = 1.068845573157468
$ws.Range("K24").Value = 1.063868455997489
$ws.Range("L24").Value = 1.070290780801616
$ws.Range("M24").Value = 1.077767639004089
$ws.Range("N24").Value = 1.070363456469652

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.065089221191928
$ws.Range("D25").Value = 1.061977647597132
$ws.Range("E25").Value = 1.06895518882095
$ws.Range("F25").Value = 1.076408968784166
$ws.Range("I25").Value = 1.050408819056365
$ws.Range("J25").Value = 1.070511256494679
$ws.Range("K25").Value = 1.064950118420838
$ws.Range("L25").Value = 1.071907040980082
$ws.Range("M25").Value = 1.079339116050998
$ws.Range("N25").Value = 1.072031505268259
